# productslist.xlsx update:
#  - Sheet1: Chiaro / Coconut Flavour Over Ice valid quantities bumped (20->120, 40->110)
#  - Sheet1: selection moved to B4, no longer the "active" (front) tab
#  - New "Sheet2" worksheet added after Sheet1: Machines names / Valid Quantity,
#    with every machine currently at a flat "1" (handles "unavailable product")
#  - Sheet2 becomes the active tab with A2 selected

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: bump the two quantities that changed ---
$ws1.Range("B2").Value = "120"
$ws1.Range("B3").Value = "110"

# Move Sheet1's own selection to B4 (it stays behind, not the active tab once
# Sheet2 is added/activated below).
$ws1.Range("B4").Select() | Out-Null

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

# Header row (row 1) - same look as Sheet1's title row: 16pt, text format.
$ws2.Range("A1:B1").NumberFormat = "@"
$ws2.Range("A1:B1").Font.Size = 16
$ws2.Range("A1").Value = "Machines names"
$ws2.Range("B1").Value = "Valid Quantity"
$ws2.Rows.Item(1).RowHeight = 19.7

# Machine names (col A) - text format, wrapped.
$machines = @(
    "Inissia",
    "Nespresso Atelier",
    "Creatista Pro",
    "Citiz",
    "Citiz Platinum",
    "Lattissima One",
    "Pixie",
    "Gran Lattissima"
)

$ws2.Range("A2:A9").NumberFormat = "@"
$ws2.Range("A2:A9").WrapText = $true

for ($i = 0; $i -lt $machines.Length; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $machines[$i]
}

# Valid quantity column (col B) - text format, flat "1" for every machine.
$ws2.Range("B2:B9").NumberFormat = "@"
for ($row = 2; $row -le 9; $row++) {
    $ws2.Range("B$row").Value = "1"
}

# Column widths roughly matching the authored sheet.
$ws2.Columns.Item(1).ColumnWidth = 24.45
$ws2.Columns.Item(2).ColumnWidth = 27.5

# Sheet2 becomes the active/front sheet, with A2 selected.
$ws2.Range("A2").Select() | Out-Null
